$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (Item ID, Item Name) for rows 7-10
$newRows = @(
    @{ Row = 7;  Id = 2876884; Name = "ISIS LS Gel 2,5L Lemon" },
    @{ Row = 8;  Id = 2917788; Name = "Le Chat LS 2,5kg bag Regular" },
    @{ Row = 9;  Id = 2918203; Name = "ISIS LS bag 2,5Kg Citron Limitless" },
    @{ Row = 10; Id = 2922764; Name = "Le Chat Reg 2,5L FRESCO" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Id
    $ws.Cells.Item($r.Row, 2).Value = $r.Name

    # Match the existing bordered style used by the prior rows (A2:D6)
    $ws.Range("A" + $r.Row + ":D" + $r.Row).Style = $ws.Range("A6:D6").Style
    $ws.Cells.Item($r.Row, 1).Borders.LineStyle = $ws.Cells.Item(6, 1).Borders.LineStyle
    $ws.Cells.Item($r.Row, 2).Borders.LineStyle = $ws.Cells.Item(6, 2).Borders.LineStyle
    $ws.Cells.Item($r.Row, 3).Borders.LineStyle = $ws.Cells.Item(6, 3).Borders.LineStyle
    $ws.Cells.Item($r.Row, 4).Borders.LineStyle = $ws.Cells.Item(6, 4).Borders.LineStyle
}

# Update selection to match final state: active cell A2, selection A2:B10
$ws.Range("A2:B10").Select()
$ws.Application.ActiveCell = $ws.Range("A2")
